$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from A4 into the new A5 label cell, matching the
# bold/bordered/centered "year" style used by A2:A4.
$ws.Cells.Item(4, 1).Copy()
$ws.Cells.Item(5, 1).PasteSpecial(-4122)

# Row 5 = 2021 data.
$ws.Cells.Item(5, 1).Value = "2021年"
$ws.Cells.Item(5, 2).Value = -58.9
$ws.Cells.Item(5, 3).Value = 24
# D5 left blank (source row has an empty cell here)
$ws.Cells.Item(5, 5).Value = 15.3
$ws.Cells.Item(5, 6).Value = -33.8
# G5 left blank (source row has an empty cell here)
$ws.Cells.Item(5, 8).Value = 17.4
$ws.Cells.Item(5, 9).Value = -20.8
$ws.Cells.Item(5, 10).Value = -17.2
$ws.Cells.Item(5, 11).Value = 89.3
$ws.Cells.Item(5, 12).Value = -77.09999999999999
$ws.Cells.Item(5, 13).Value = 11
$ws.Cells.Item(5, 14).Value = -11
$ws.Cells.Item(5, 15).Value = -3.7
$ws.Cells.Item(5, 16).Value = -37.3
$ws.Cells.Item(5, 17).Value = -61.6
$ws.Cells.Item(5, 18).Value = -86.5
# S5 left blank (source row has an empty cell here)
$ws.Cells.Item(5, 20).Value = 40.5
$ws.Cells.Item(5, 21).Value = 8.699999999999999
$ws.Cells.Item(5, 22).Value = -3.5
$ws.Cells.Item(5, 23).Value = 25.4
$ws.Cells.Item(5, 24).Value = 51.4
# Y5 left blank (source row has an empty cell here)
# Z5 left blank (source row has an empty cell here)
$ws.Cells.Item(5, 27).Value = -30.4
$ws.Cells.Item(5, 28).Value = -33.1
$ws.Cells.Item(5, 29).Value = -40.3
$ws.Cells.Item(5, 30).Value = -73.3
$ws.Cells.Item(5, 31).Value = 4.8
$ws.Cells.Item(5, 32).Value = -34.3
# AG5 left blank (source row has an empty cell here)
$ws.Cells.Item(5, 34).Value = 11.3
$ws.Cells.Item(5, 35).Value = 82.90000000000001
$ws.Cells.Item(5, 36).Value = 5.9
$ws.Cells.Item(5, 37).Value = -68.2
$ws.Cells.Item(5, 38).Value = -16
$ws.Cells.Item(5, 39).Value = 280.4
$ws.Cells.Item(5, 40).Value = -29.5
$ws.Cells.Item(5, 41).Value = -28.2
$ws.Cells.Item(5, 42).Value = -94.40000000000001
$ws.Cells.Item(5, 43).Value = -42.9
$ws.Cells.Item(5, 44).Value = 97.7
# AS5 left blank (source row has an empty cell here)
$ws.Cells.Item(5, 46).Value = -24
$ws.Cells.Item(5, 47).Value = 127.6
$ws.Cells.Item(5, 48).Value = 3.5
$ws.Cells.Item(5, 49).Value = 79.40000000000001
$ws.Cells.Item(5, 50).Value = -13.4
$ws.Cells.Item(5, 51).Value = -0.7
$ws.Cells.Item(5, 52).Value = 6.9
$ws.Cells.Item(5, 53).Value = -9.1
$ws.Cells.Item(5, 54).Value = -26.6
$ws.Cells.Item(5, 55).Value = 4.7
# BD5 left blank (source row has an empty cell here)
$ws.Cells.Item(5, 57).Value = 62.4
$ws.Cells.Item(5, 58).Value = 0.9
$ws.Cells.Item(5, 59).Value = 22.2
$ws.Cells.Item(5, 60).Value = -30.4
$ws.Cells.Item(5, 61).Value = -2.1
$ws.Cells.Item(5, 62).Value = 23.8
$ws.Cells.Item(5, 63).Value = -53.2
$ws.Cells.Item(5, 64).Value = -1.2
$ws.Cells.Item(5, 65).Value = 25.8
$ws.Cells.Item(5, 66).Value = -23.5
$ws.Cells.Item(5, 67).Value = -47.5
$ws.Cells.Item(5, 68).Value = -25.9
# BQ5 left blank (source row has an empty cell here)
$ws.Cells.Item(5, 70).Value = 53.1
$ws.Cells.Item(5, 71).Value = -5.3
$ws.Cells.Item(5, 72).Value = 10.2
$ws.Cells.Item(5, 73).Value = -37.9
$ws.Cells.Item(5, 74).Value = -34.5
$ws.Cells.Item(5, 75).Value = -44.6
$ws.Cells.Item(5, 76).Value = 76.3
$ws.Cells.Item(5, 77).Value = 34.3
$ws.Cells.Item(5, 78).Value = 48.4
$ws.Cells.Item(5, 79).Value = -20.5
$ws.Cells.Item(5, 80).Value = 5.4
$ws.Cells.Item(5, 81).Value = 30
$ws.Cells.Item(5, 82).Value = -66.90000000000001
$ws.Cells.Item(5, 83).Value = -17.4
$ws.Cells.Item(5, 84).Value = 10.4
$ws.Cells.Item(5, 85).Value = 19.4
$ws.Cells.Item(5, 86).Value = 101.3
$ws.Cells.Item(5, 87).Value = 5.5
$ws.Cells.Item(5, 88).Value = -24.3
$ws.Cells.Item(5, 89).Value = -18.2
$ws.Cells.Item(5, 90).Value = -3.5
$ws.Cells.Item(5, 91).Value = -19.5
$ws.Cells.Item(5, 92).Value = 190.9
$ws.Cells.Item(5, 93).Value = 26.6
$ws.Cells.Item(5, 94).Value = -30.6
$ws.Cells.Item(5, 95).Value = -9
# CR5 left blank (source row has an empty cell here)
$ws.Cells.Item(5, 97).Value = 28
$ws.Cells.Item(5, 98).Value = 21.8
$ws.Cells.Item(5, 99).Value = -14.8
$ws.Cells.Item(5, 100).Value = -40.8
$ws.Cells.Item(5, 101).Value = -64.5
$ws.Cells.Item(5, 102).Value = 14.2
$ws.Cells.Item(5, 103).Value = 25.2
$ws.Cells.Item(5, 104).Value = 120.5
$ws.Cells.Item(5, 105).Value = 4.9
$ws.Cells.Item(5, 106).Value = -12.8
$ws.Cells.Item(5, 107).Value = 101.6
$ws.Cells.Item(5, 108).Value = -89.59999999999999
$ws.Cells.Item(5, 109).Value = 5.5
$ws.Cells.Item(5, 110).Value = 2.6
$ws.Cells.Item(5, 111).Value = -14.1
$ws.Cells.Item(5, 112).Value = -16.1
$ws.Cells.Item(5, 113).Value = 11.3
$ws.Cells.Item(5, 114).Value = 91.59999999999999
$ws.Cells.Item(5, 115).Value = 54.9
